# Generate Report for Handoff
# Updates the "Status" text from "Handed back: in sync with en-US" to "Ready for handoff"
# and refreshes the associated timestamps, across the Overview, zh-cn, and de-de sheets.
# Also narrows the "Latest Handoff Datetime" / "Latest Handback DateTime" style columns.

$wb = $excel.ActiveWorkbook

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Ready for handoff"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-12 09:12:40"
$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-12 09:12:33"
$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-12 09:12:40"
$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
